$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("WC49 P5H", "Detección de sealling mal puesto", "2024-05-29", "12:42:27", "Mañana", "12:42:29", "0:00:02", "N/A"),
    @("WC49 P5H", "Screw K30 no lo detecta puesto", "2024-05-29", "12:42:32", "Mañana", "12:42:34", "0:00:02", "N/A"),
    @("WV50 FILTER", "Fallo cámara ferrite", "2024-05-29", "12:45:18", "Mañana", "12:45:19", "0:00:01", "N/A"),
    @("WV50 FILTER", "NOK Soldad. Plástico+Metal", "2024-05-29", "12:45:25", "Mañana", "12:45:27", "0:00:02", "N/A"),
    @("WV50 FILTER", "Fallo cámara QR", "2024-05-29", "12:45:29", "Mañana", "12:45:36", "0:00:07", "0.13 minutos"),
    @("WV50 FILTER", "Fallo cámara QR", "2024-05-29", "12:45:31", "Mañana", "12:45:32", "0:00:01", "0.10 minutos"),
    @("WC49 P5H", "Atasco tuerca", "2024-05-29", "12:55:16", "Mañana", "12:55:21", "0:00:05", "N/A"),
    @("WC49 P5H", "Tornillo atascado", "2024-05-29", "12:55:29", "Mañana", "12:55:31", "0:00:02", "N/A"),
    @("WC49 P5H", "La cámara no detecta Busbar", "2024-05-29", "12:55:35", "Mañana", "12:55:37", "0:00:02", "0.21 minutos"),
    @("WC49 P5H", "Screw K30 no lo detecta puesto", "2024-05-29", "12:56:52", "Mañana", "12:56:54", "0:00:02", "N/A"),
    @("WC49 P5H", "No coloca bien el sealling", "2024-05-29", "12:56:55", "Mañana", "12:56:57", "0:00:02", "N/A"),
    @("WC47 NACP", "Etiquetadora", "2024-05-29", "12:57:43", "Mañana", "12:57:51", "0:00:08", "N/A"),
    @("WC47 NACP", "Fallo tornillo", "2024-05-29", "12:57:47", "Mañana", "12:57:52", "0:00:05", "N/A"),
    @("WC47 NACP", "Fallo cámara visión", "2024-05-29", "12:57:50", "Mañana", "12:57:52", "0:00:02", "0.06 minutos"),
    @("WC47 NACP", "AOI no detecta pieza", "2024-05-29", "12:57:54", "Mañana", "12:57:56", "0:00:02", "0.06 minutos")
)

$startRow = 16
$ws.Range("C16:C30").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
